$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet1 ("Metadata") ---

# Row 3: Version
$ws.Range("B3").Value = "0.1.7"

# Row 6: Status
$ws.Range("B6").Value = "draft"

# Row 8: Date
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Row 10: Contact (existing row updated to new publisher contact text)
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11 already exists as a duplicate "Contact" row - give it the second
# contact entry (Bob Milius).
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Make room for a new "Jurisdiction" row at position 12 by shifting the
# remaining rows (old Description/Purpose/Copyright/Immutable, rows 12-15)
# down to rows 13-16. Work bottom-up so we don't clobber source rows before
# they're copied. First extend formatting to the new row 16 using the
# format of row 15 (the last existing row) so the new row matches the rest
# of the table, then copy values into place.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($r = 15; $r -ge 12; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Value()
}

# Row 12: new Jurisdiction property (empty value)
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
